$d = $word.ActiveDocument

# --- 1. Remove the "Задайте най-подходящите типове данни... Попълнете всяка
#        таблица с точно 3 записа..." paragraph that used to precede the
#        "База данни Hotel" (Heading2) paragraph, and move the
#        <w:lastRenderedPageBreak/> marker onto the first run of that heading
#        paragraph. (The target diff deletes the whole paragraph and the
#        page-break hint now opens the heading run instead.)
#        NOTE: there is a near-identical paragraph elsewhere in the document
#        that talks about "точно 5 записа" instead of "точно 3" - match on
#        the "точно 3" marker specifically so that one is left untouched.

$targetIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    $t = $p.Range.Text
    if ($t -like "*точно 3*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 0) {
    $target = $d.Paragraphs($targetIdx)
    $target.Range.Delete()

    # Re-fetch the (now shifted) heading paragraph fresh from the collection
    # so we don't operate on a stale Range left over from before the delete.
    $heading = $d.Paragraphs($targetIdx)
    $headingRange = $heading.Range

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="60A3E55B" w14:textId="77777777" w:rsidR="009365EE" w:rsidRPr="009365EE" w:rsidRDefault="009365EE" w:rsidP="001C681C"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r w:rsidRPr="009365EE"><w:lastRenderedPageBreak/><w:t xml:space="preserve">База данни </w:t></w:r><w:r w:rsidRPr="009365EE"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hotel</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $headingRange.InsertXML($xml)
}

# --- 2. Remove the trailing "Задайте най-подходящите типове данни... Попълнете
#        всяка таблица само с 3 записа..." paragraph at the very end of the
#        document (right before the section properties). It is simply
#        deleted, nothing follows it.

$trailingIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    $t = $p.Range.Text
    if ($t -like "*само с*") {
        $trailingIdx = $i
        break
    }
}

if ($trailingIdx -gt 0) {
    $trailing = $d.Paragraphs($trailingIdx)
    $trailing.Range.Delete()
}
